# DataTable.xlsx update: refresh the "Highest qualification level by age
# and gender" (APS) row with the newer NCF-recoded data period, replacing
# the old caveat text that said the period/date was still TBC.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = "Jan 2023 - Dec 2023 (16/04/24)"
$ws.Range("D7").Value = "Apr 2023 - Mar 2024 (18/07/24)"

# Reflect where the author was last working when they saved: scrolled so
# row 2 is at the top of the viewport, with C7:D7 (the cells just edited)
# selected.
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C7:D7").Select()
